$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rename the "dateCreated" header (D1) to "date" - this is the core content edit.
# (Excel's save process will naturally renumber/compact the shared-string table,
# which is what produces the rest of the <v> index churn seen in the diff.)
$ws.Range("D1").Value = "date"

# Selection moves to E1 (single cell) and the view scrolls back to show column A.
$ws.Range("E1").Select()

# Row heights grow (content/font metrics changed on resave).
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(3).RowHeight = 105
